{"js": "const pairs = [\n  [\"999\u00f75=199, 4\", \"733\u00f74=183, 1\"],\n  [\"672\u00f72=336, 0\", \"160\u00f78=20, 0\"],\n  [\"433\u00f72=216, 1\", \"214\u00f75=42, 4\"],\n  [\"670\u00f77=95, 5\", \"408\u00f79=45, 3\"],\n  [\"229\u00f77=32, 5\", \"355\u00f79=39, 4\"],\n  [\"655\u00f73=218, 1\", \"655\u00f79=72, 7\"],\n  [\"730\u00f74=182, 2\", \"928\u00f72=464, 0\"],\n  [\"176\u00f75=35, 1\", \"987\u00f77=141, 0\"],\n  [\"690\u00f75=138, 0\", \"288\u00f78=36, 0\"],\n  [\"834\u00f77=119, 1\", \"250\u00f76=41, 4\"],\n  [\"539\u00f72=269, 1\", \"607\u00f78=75, 7\"],\n  [\"956\u00f76=159, 2\", \"566\u00f75=113, 1\"],\n  [\"940\u00f79=104, 4\", \"243\u00f73=81, 0\"],\n  [\"480\u00f78=60, 0\", \"232\u00f75=46, 2\"],\n  [\"738\u00f74=184, 2\", \"457\u00f79=50, 7\"],\n  [\"591\u00f77=84, 3\", \"791\u00f76=131, 5\"],\n  [\"673\u00f76=112, 1\", \"316\u00f78=39, 4\"],\n  [\"229\u00f72=114, 1\", \"948\u00f73=316, 0\"],\n  [\"741\u00f79=82, 3\", \"804\u00f76=134, 0\"],\n  [\"301\u00f78=37, 5\", \"825\u00f78=103, 1\"],\n  [\"911\u00f77=130, 1\", \"881\u00f73=293, 2\"],\n  [\"265\u00f72=132, 1\", \"910\u00f77=130, 0\"],\n  [\"297\u00f72=148, 1\", \"726\u00f75=145, 1\"],\n  [\"533\u00f72=266, 1\", \"522\u00f73=174, 0\"],\n  [\"674\u00f79=74, 8\", \"654\u00f77=93, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"999\u00f75=199, 4\", \"733\u00f74=183, 1\"),\n    @(\"672\u00f72=336, 0\", \"160\u00f78=20, 0\"),\n    @(\"433\u00f72=216, 1\", \"214\u00f75=42, 4\"),\n    @(\"670\u00f77=95, 5\", \"408\u00f79=45, 3\"),\n    @(\"229\u00f77=32, 5\", \"355\u00f79=39, 4\"),\n    @(\"655\u00f73=218, 1\", \"655\u00f79=72, 7\"),\n    @(\"730\u00f74=182, 2\", \"928\u00f72=464, 0\"),\n    @(\"176\u00f75=35, 1\", \"987\u00f77=141, 0\"),\n    @(\"690\u00f75=138, 0\", \"288\u00f78=36, 0\"),\n    @(\"834\u00f77=119, 1\", \"250\u00f76=41, 4\"),\n    @(\"539\u00f72=269, 1\", \"607\u00f78=75, 7\"),\n    @(\"956\u00f76=159, 2\", \"566\u00f75=113, 1\"),\n    @(\"940\u00f79=104, 4\", \"243\u00f73=81, 0\"),\n    @(\"480\u00f78=60, 0\", \"232\u00f75=46, 2\"),\n    @(\"738\u00f74=184, 2\", \"457\u00f79=50, 7\"),\n    @(\"591\u00f77=84, 3\", \"791\u00f76=131, 5\"),\n    @(\"673\u00f76=112, 1\", \"316\u00f78=39, 4\"),\n    @(\"229\u00f72=114, 1\", \"948\u00f73=316, 0\"),\n    @(\"741\u00f79=82, 3\", \"804\u00f76=134, 0\"),\n    @(\"301\u00f78=37, 5\", \"825\u00f78=103, 1\"),\n    @(\"911\u00f77=130, 1\", \"881\u00f73=293, 2\"),\n    @(\"265\u00f72=132, 1\", \"910\u00f77=130, 0\"),\n    @(\"297\u00f72=148, 1\", \"726\u00f75=145, 1\"),\n    @(\"533\u00f72=266, 1\", \"522\u00f73=174, 0\"),\n    @(\"674\u00f79=74, 8\", \"654\u00f77=93, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
